# Apply the "Request; Worker; Server folder contains the latest working
# code." edit to the Password Cracking Stats workbook.
#
# The underlying change is simple: fill in previously-missing raw trial
# data on the "Stat Trial 1" sheet (rows 9/10 and 19/20, columns K:O plus
# D10/D20). Every other value that differs in the diff (rows 29/30/33/39/40
# and the cached chart points) is a pure formula/chart-cache recalculation
# that falls out of these raw inputs automatically.
#
# It also nudges the workbook window width and the active sheet's
# scroll/selection state to match the author's saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stat Trial 1")

# --- Row 9 / Row 19 (length-10..14 raw values) ---------------------------
$ws.Range("K9").Value  = 155156
$ws.Range("L9").Value  = 131061
$ws.Range("M9").Value  = 457032
$ws.Range("N9").Value  = 483453
$ws.Range("O9").Value  = 980804

$ws.Range("K19").Value = 2564
$ws.Range("L19").Value = 3063
$ws.Range("M19").Value = 1275
$ws.Range("N19").Value = 846
$ws.Range("O19").Value = 1694

# --- Row 10 / Row 20 (length-3 raw values) --------------------------------
$ws.Range("D10").Value = 995654
$ws.Range("D20").Value = 3636

# --- Workbook window + sheet view/selection state -------------------------
$excel.ActiveWindow.Width = 7110

$ws.Activate()
$ws.Application.Goto($ws.Range("K1"), $true)
$ws.Range("P9").Select()
